# Refresh the regression-results table (category / Coef. / Std.Err. / t / P>|t| / [0.025 / 0.975] / coef_pos)
# with the re-run hourly dataset. Bucket breakdown near the bottom of the gen/cap
# distribution is now finer (5%/10%/15%/20% instead of a single "20%>gen/cap" bucket),
# so the table grows from 25 data rows to 28 (A1:H26 -> A1:H29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = 28
$colCount = 8
$data = New-Object "object[,]" $rowCount,$colCount

# row 2: 10%>gen/cap>5%
$data[0,0] = '10%>gen/cap>5%'
$data[0,1] = [double]"0.0"
$data[0,2] = ""
$data[0,3] = ""
$data[0,4] = ""
$data[0,5] = ""
$data[0,6] = ""
$data[0,7] = [double]"0.1880725913802447"

# row 3: 100%>gen/cap>95%
$data[1,0] = '100%>gen/cap>95%'
$data[1,1] = [double]"0.09559710390669388"
$data[1,2] = ""
$data[1,3] = ""
$data[1,4] = ""
$data[1,5] = ""
$data[1,6] = ""
$data[1,7] = [double]"0.2836696952869385"

# row 4: 105%>gen/cap>100%
$data[2,0] = '105%>gen/cap>100%'
$data[2,1] = [double]"0.1510524751983966"
$data[2,2] = [double]"0.01297507593763596"
$data[2,3] = [double]"17.18719027254989"
$data[2,4] = [double]"0.05771751675679879"
$data[2,5] = [double]"0.1256178560764017"
$data[2,6] = [double]"0.1764870943203905"
$data[2,7] = [double]"0.3391250665786413"

# row 5: 110%>gen/cap>105%
$data[3,0] = '110%>gen/cap>105%'
$data[3,1] = [double]"0.1779845313118368"
$data[3,2] = [double]"0.01359658199882654"
$data[3,3] = [double]"20.57171661340502"
$data[3,4] = [double]"0.01580215701252231"
$data[3,5] = [double]"0.1513317968817242"
$data[3,6] = [double]"0.2046372657419491"
$data[3,7] = [double]"0.3660571226920814"

# row 6: 115%>gen/cap>110%
$data[4,0] = '115%>gen/cap>110%'
$data[4,1] = [double]"0.04453413997050377"
$data[4,2] = [double]"0.006668347637500935"
$data[4,3] = [double]"4.079605483349251"
$data[4,4] = [double]"0.01681100044957815"
$data[4,5] = [double]"0.0314624787243593"
$data[4,6] = [double]"0.05760580121664776"
$data[4,7] = [double]"0.2326067313507484"

# row 7: 120%>gen/cap>115%
$data[5,0] = '120%>gen/cap>115%'
$data[5,1] = [double]"0.01043783477834273"
$data[5,2] = [double]"0.00880152424180306"
$data[5,3] = [double]"0.5323332147068555"
$data[5,4] = [double]"0.03887219234474882"
$data[5,5] = [double]"-0.006815331702531453"
$data[5,6] = [double]"0.0276910012592169"
$data[5,7] = [double]"0.1985104261585874"

# row 8: 125%>gen/cap>120%
$data[6,0] = '125%>gen/cap>120%'
$data[6,1] = [double]"0.003933398020328643"
$data[6,2] = [double]"0.00020238629662523"
$data[6,3] = [double]"0.2837240952940529"
$data[6,4] = [double]"7.927817321965518e-83"
$data[6,5] = [double]"0.003536626055141385"
$data[6,6] = [double]"0.0043301699855159"
$data[6,7] = [double]"0.1920059894005733"

# row 9: 130%>gen/cap>125%
$data[7,0] = '130%>gen/cap>125%'
$data[7,1] = [double]"0.004015826779403169"
$data[7,2] = [double]"0.0002499802908057933"
$data[7,3] = [double]"0.2345193225812576"
$data[7,4] = [double]"2.062997279377916e-58"
$data[7,5] = [double]"0.003525748286466652"
$data[7,6] = [double]"0.004505905272339695"
$data[7,7] = [double]"0.1920884181596478"

# row 10: 135%>gen/cap>130%
$data[8,0] = '135%>gen/cap>130%'
$data[8,1] = [double]"0.004223004832812028"
$data[8,2] = [double]"0.0002380450526303365"
$data[8,3] = [double]"0.2589833517112979"
$data[8,4] = [double]"4.773060188465943e-70"
$data[8,5] = [double]"0.00375632499870182"
$data[8,6] = [double]"0.004689684666922237"
$data[8,7] = [double]"0.1922955962130567"

# row 11: 15%>gen/cap>10%
$data[9,0] = '15%>gen/cap>10%'
$data[9,1] = [double]"0.03148956516287641"
$data[9,2] = ""
$data[9,3] = ""
$data[9,4] = ""
$data[9,5] = ""
$data[9,6] = ""
$data[9,7] = [double]"0.2195621565431211"

# row 12: 20%>gen/cap>15%
$data[10,0] = '20%>gen/cap>15%'
$data[10,1] = [double]"0.04009389428101376"
$data[10,2] = ""
$data[10,3] = ""
$data[10,4] = ""
$data[10,5] = ""
$data[10,6] = ""
$data[10,7] = [double]"0.2281664856612584"

# row 13: 25%>gen/cap>20%
$data[11,0] = '25%>gen/cap>20%'
$data[11,1] = [double]"0.0476860042510685"
$data[11,2] = ""
$data[11,3] = ""
$data[11,4] = ""
$data[11,5] = ""
$data[11,6] = ""
$data[11,7] = [double]"0.2357585956313132"

# row 14: 30%>gen/cap>25%
$data[12,0] = '30%>gen/cap>25%'
$data[12,1] = [double]"0.056054627564771"
$data[12,2] = [double]"0.009040940971138788"
$data[12,3] = [double]"8.126623123894474"
$data[12,4] = [double]"0.04980304778512183"
$data[12,5] = [double]"0.03833170627615557"
$data[12,6] = [double]"0.07377754885338629"
$data[12,7] = [double]"0.2441272189450157"

# row 15: 35%>gen/cap>30%
$data[13,0] = '35%>gen/cap>30%'
$data[13,1] = [double]"0.05949143801937335"
$data[13,2] = [double]"0.008845465359062419"
$data[13,3] = [double]"9.150753779538713"
$data[13,4] = [double]"0.05251808944827788"
$data[13,5] = [double]"0.04215170068035005"
$data[13,6] = [double]"0.07683117535839665"
$data[13,7] = [double]"0.247564029399618"

# row 16: 40%>gen/cap>35%
$data[14,0] = '40%>gen/cap>35%'
$data[14,1] = [double]"0.06285931152351278"
$data[14,2] = [double]"0.008830606984105023"
$data[14,3] = [double]"9.281126380735753"
$data[14,4] = [double]"0.03742866811173217"
$data[14,5] = [double]"0.04554865558109091"
$data[14,6] = [double]"0.0801699674659346"
$data[14,7] = [double]"0.2509319029037574"

# row 17: 45%>gen/cap>40%
$data[15,0] = '45%>gen/cap>40%'
$data[15,1] = [double]"0.06568699870112284"
$data[15,2] = [double]"0.00868922225373481"
$data[15,3] = [double]"10.15076775283177"
$data[15,4] = [double]"0.05265612293179296"
$data[15,5] = [double]"0.04865359614103339"
$data[15,6] = [double]"0.08272040126121234"
$data[15,7] = [double]"0.2537595900813675"

# row 18: 5%>gen/cap
$data[16,0] = '5%>gen/cap'
$data[16,1] = [double]"-0.1880725913802447"
$data[16,2] = [double]"0.01170277594891895"
$data[16,3] = [double]"-33.28944777559816"
$data[16,4] = [double]"0.01467037756862944"
$data[16,5] = [double]"-0.2110134259768802"
$data[16,6] = [double]"-0.1651317567836091"
$data[16,7] = [double]"0.0"

# row 19: 50%>gen/cap>45%
$data[17,0] = '50%>gen/cap>45%'
$data[17,1] = [double]"0.06599946126281732"
$data[17,2] = [double]"0.008809042061937472"
$data[17,3] = [double]"9.902846044563343"
$data[17,4] = [double]"0.04741340693143789"
$data[17,5] = [double]"0.0487312023650459"
$data[17,6] = [double]"0.08326772016058885"
$data[17,7] = [double]"0.254072052643062"

# row 20: 55%>gen/cap>50%
$data[18,0] = '55%>gen/cap>50%'
$data[18,1] = [double]"0.06567866701675514"
$data[18,2] = [double]"0.009205147573644205"
$data[18,3] = [double]"9.231897571191439"
$data[18,4] = [double]"0.04759618860331526"
$data[18,5] = [double]"0.04763398017173601"
$data[18,6] = [double]"0.0837233538617742"
$data[18,7] = [double]"0.2537512583969998"

# row 21: 60%>gen/cap>55%
$data[19,0] = '60%>gen/cap>55%'
$data[19,1] = [double]"0.0663915239330164"
$data[19,2] = [double]"0.009438816938096518"
$data[19,3] = [double]"9.433950381713014"
$data[19,4] = [double]"0.05164104155379814"
$data[19,5] = [double]"0.04788874738821575"
$data[19,6] = [double]"0.08489430047781708"
$data[19,7] = [double]"0.2544641153132611"

# row 22: 65%>gen/cap>60%
$data[20,0] = '65%>gen/cap>60%'
$data[20,1] = [double]"0.06899438716721364"
$data[20,2] = [double]"0.009417296868933465"
$data[20,3] = [double]"9.861490071869861"
$data[20,4] = [double]"0.0479608699813569"
$data[20,5] = [double]"0.05053380126127676"
$data[20,6] = [double]"0.08745497307315046"
$data[20,7] = [double]"0.2570669785474583"

# row 23: 70%>gen/cap>65%
$data[21,0] = '70%>gen/cap>65%'
$data[21,1] = [double]"0.0708686706739472"
$data[21,2] = [double]"0.009560548361665559"
$data[21,3] = [double]"9.507439581667656"
$data[21,4] = [double]"0.05735888357781408"
$data[21,5] = [double]"0.05212728609443056"
$data[21,6] = [double]"0.0896100552534638"
$data[21,7] = [double]"0.2589412620541919"

# row 24: 75%>gen/cap>70%
$data[22,0] = '75%>gen/cap>70%'
$data[22,1] = [double]"0.07275187338324948"
$data[22,2] = [double]"0.00995993160667035"
$data[22,3] = [double]"9.885765084249979"
$data[22,4] = [double]"0.06492683325242563"
$data[22,5] = [double]"0.0532275749294217"
$data[22,6] = [double]"0.09227617183707726"
$data[22,7] = [double]"0.2608244647634941"

# row 25: 80%>gen/cap>75%
$data[23,0] = '80%>gen/cap>75%'
$data[23,1] = [double]"0.08063055286310358"
$data[23,2] = [double]"0.01027781763553838"
$data[23,3] = [double]"10.89631900024567"
$data[23,4] = [double]"0.05915340984844311"
$data[23,5] = [double]"0.06048310811986542"
$data[23,6] = [double]"0.1007779976063416"
$data[23,7] = [double]"0.2687031442433482"

# row 26: 85%>gen/cap>80%
$data[24,0] = '85%>gen/cap>80%'
$data[24,1] = [double]"0.0845418936153209"
$data[24,2] = [double]"0.01032596987816424"
$data[24,3] = [double]"11.4840035256946"
$data[24,4] = [double]"0.06881748857141784"
$data[24,5] = [double]"0.06430008386140391"
$data[24,6] = [double]"0.1047837033692379"
$data[24,7] = [double]"0.2726144849955656"

# row 27: 90%>gen/cap>85%
$data[25,0] = '90%>gen/cap>85%'
$data[25,1] = [double]"0.08215563187890858"
$data[25,2] = [double]"0.01042225613604837"
$data[25,3] = [double]"9.984956131143194"
$data[25,4] = [double]"0.07805650701418508"
$data[25,5] = [double]"0.06172515760592509"
$data[25,6] = [double]"0.1025861061518919"
$data[25,7] = [double]"0.2702282232591532"

# row 28: 95%>gen/cap>90%
$data[26,0] = '95%>gen/cap>90%'
$data[26,1] = [double]"0.08248859821366536"
$data[26,2] = [double]"0.01128897907392203"
$data[26,3] = [double]"9.512242140003572"
$data[26,4] = [double]"0.1045537621422679"
$data[26,5] = [double]"0.06035896668564295"
$data[26,6] = [double]"0.1046182297416878"
$data[26,7] = [double]"0.27056118959391"

# row 29: gen/cap>135%
$data[27,0] = 'gen/cap>135%'
$data[27,1] = [double]"0.004203108549082341"
$data[27,2] = [double]"0.0001787439829560479"
$data[27,3] = [double]"0.3432800806883717"
$data[27,4] = [double]"1.579459369176043e-117"
$data[27,5] = [double]"0.003852686595790732"
$data[27,6] = [double]"0.004553530502373951"
$data[27,7] = [double]"0.192275699929327"

$ws.Range("A2:H29").Value = $data

